$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty "risk of bias" cells in column I with "Low"
# (matching the value already used in the same rows' H/J columns), for the
# dose-response/rob rows added in this pass. For the rows whose I cell needs
# to pick up the "always-4" shading used by column C in that row (9, 11, 13,
# 15), copy C's cell format into I first so the style index lines up with the
# rest of that shaded column; the other rows keep their existing format.

$rowsNeedingFormatFix = @(9, 11, 13, 15)
foreach ($r in $rowsNeedingFormatFix) {
    $ws.Range("C$r").Copy()
    $ws.Range("I$r").PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = 0

$rowsToFill = @(8, 9, 10, 11, 12, 13, 14, 15, 16, 78, 79, 80)
foreach ($r in $rowsToFill) {
    $ws.Range("I$r").Value = "Low"
}

# Restore the saved selection to match the author's last interaction: the
# entire row 16 selected with A16 as the active cell.
$ws.Rows.Item(16).Select()
